# Use 6.9% DM not 5.1% for digestate
$wb = $excel.ActiveWorkbook

$slurry = $wb.Worksheets.Item("Slurry")

# Digestate ("Afgasset biomasse") dry-matter rows: update 5.1 -> 6.9
$slurry.Range("C4").Value = 6.9
$slurry.Range("C7").Value = 6.9

# Make "Slurry" the active sheet/tab, with C8 selected (also clears the
# previously-active "Climate" tab's tabSelected flag).
$slurry.Activate()
$slurry.Range("C8").Select()
